# Update countries & provincias Spain
#
# The underlying sharedStrings order shuffled three country names
# (Estado de Palestina, Isla de Man, Guayana Francesa) up the list, which
# shows up as new country labels landing on rows 106-109, 114-115 and
# 136-141 (the row's numeric stats are untouched by the rename - they are
# updated separately below to match the refreshed case counts).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update country names (reordering of Estado de Palestina, Isla de Man, Guayana Francesa)
$ws.Range("A106").Value = "Estado de Palestina"
$ws.Range("A107").Value = "Jordania"
$ws.Range("A108").Value = "Taiwan"
$ws.Range("A109").Value = "Reunion"
$ws.Range("A114").Value = "Isla de Man"
$ws.Range("A115").Value = "Consejo Danes para los Refugiados"
$ws.Range("A136").Value = "Guayana Francesa"
$ws.Range("A137").Value = "Gabon"
$ws.Range("A138").Value = "Aruba"
$ws.Range("A139").Value = "Tanzania"
$ws.Range("A140").Value = "Monaco"
$ws.Range("A141").Value = "Etiopia"

# Update statistics values
$ws.Range("E18").Value = 9544
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 1288
$ws.Range("E34").Value = 6716
$ws.Range("F34").Value = 63
$ws.Range("G34").Value = 5
$ws.Range("H34").Value = 157
$ws.Range("B35").Value = 6523
$ws.Range("C35").Value = 55
$ws.Range("D35").Value = 3819
$ws.Range("E35").Value = 2639
$ws.Range("F35").Value = 60
$ws.Range("G35").Value = 2
$ws.Range("H35").Value = 65
$ws.Range("B43").Value = 5251
$ws.Range("C43").Value = 69
$ws.Range("D43").Value = 2967
$ws.Range("E43").Value = 2198
$ws.Range("F43").Value = 51
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = 86
$ws.Range("C106").Value = 28
$ws.Range("D106").Value = 69
$ws.Range("E106").Value = 331
$ws.Range("F106").Value = 0
$ws.Range("H106").Value = 2
$ws.Range("B107").Value = 402
$ws.Range("D107").Value = 259
$ws.Range("E107").Value = 136
$ws.Range("F107").Value = 5
$ws.Range("H107").Value = 7
$ws.Range("B108").Value = 395
$ws.Range("D108").Value = 166
$ws.Range("E108").Value = 223
$ws.Range("F108").Value = 0
$ws.Range("H108").Value = 6
$ws.Range("B109").Value = 394
$ws.Range("D109").Value = 237
$ws.Range("E109").Value = 157
$ws.Range("F109").Value = 4
$ws.Range("H109").Value = 0
$ws.Range("D110").Value = 79
$ws.Range("E110").Value = 288
$ws.Range("B114").Value = 289
$ws.Range("C114").Value = 5
$ws.Range("D114").Value = 168
$ws.Range("E114").Value = 117
$ws.Range("F114").Value = 10
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 4
$ws.Range("B115").Value = 287
$ws.Range("C115").Value = 20
$ws.Range("D115").Value = 25
$ws.Range("E115").Value = 239
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 23
$ws.Range("D118").Value = 70
$ws.Range("E118").Value = 161
$ws.Range("D130").Value = 112
$ws.Range("E130").Value = 23
$ws.Range("B136").Value = 96
$ws.Range("C136").Value = 10
$ws.Range("D136").Value = 61
$ws.Range("E136").Value = 35
$ws.Range("F136").Value = 2
$ws.Range("H136").Value = 0
$ws.Range("C137").Value = 15
$ws.Range("D137").Value = 6
$ws.Range("E137").Value = 88
$ws.Range("F137").Value = 0
$ws.Range("H137").Value = 1
$ws.Range("B138").Value = 95
$ws.Range("D138").Value = 39
$ws.Range("E138").Value = 54
$ws.Range("F138").Value = 1
$ws.Range("H138").Value = 2
$ws.Range("B139").Value = 94
$ws.Range("D139").Value = 11
$ws.Range("E139").Value = 79
$ws.Range("F139").Value = 0
$ws.Range("H139").Value = 4
$ws.Range("B140").Value = 93
$ws.Range("D140").Value = 12
$ws.Range("E140").Value = 78
$ws.Range("F140").Value = 2
$ws.Range("B141").Value = 92
$ws.Range("D141").Value = 15
$ws.Range("E141").Value = 74
$ws.Range("H141").Value = 3